# Commit: "added slide titles for all slides"
#
# The deck starts with a single title slide ("The UNIX CLI"). This
# script renames that title and appends the ten follow-on session
# slides (each with just a title filled in on a "Title and Content"
# layout, as the commit only populated titles), reproducing the final
# slide order from the presentation.xml <p:sldIdLst>. Slide 3 keeps its
# bullet list ("UNIX-like operating systems") and slide 2 gets its
# speaker note, matching the target OOXML.

$p = $ppt.ActivePresentation

# 1. Retitle the existing title slide.
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "The UNIX command line"

# 2. Append the new content slides (Title and Content layout = 2).
$p.Slides.Add(2, 2)
$p.Slides.Add(3, 2)
$p.Slides.Add(4, 2)
$p.Slides.Add(5, 2)
$p.Slides.Add(6, 2)
$p.Slides.Add(7, 2)
$p.Slides.Add(8, 2)
$p.Slides.Add(9, 2)
$p.Slides.Add(10, 2)
$p.Slides.Add(11, 2)

$titles = @(
    "History of UNIX",
    "UNIX-like operating systems",
    "UNIX file structure conventions",
    "UNIX pipes",
    "Environment variables",
    "File descriptors",
    "File extensions",
    "Commonly-used CLI commands",
    "Compiling",
    "Exercise: UNIX commands"
)

for ($i = 0; $i -lt $titles.Count; $i++) {
    $slide = $p.Slides.Item($i + 2)
    $slide.Shapes.Item(1).TextFrame.TextRange.Text = $titles[$i]
}

# 3. Slide 3 ("UNIX-like operating systems") carries real bullet content.
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2).TextFrame.TextRange
$body3.Text = "Mac OSX`rLinux`rCygwin`rVirtualization"

# 4. Slide 2 ("History of UNIX") gets its speaker note.
$notes2 = $p.Slides.Item(2).NotesPage
$notesBody2 = $notes2.Shapes.Placeholders.Item(2)
$notesBody2.TextFrame.TextRange.Text = "Re-use ComPhy exercises"
